$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new drug (well 10 / column K) to the plate layout metadata sheets:
#   row A (row 2) -> Temozolomide
#   row B (row 3) -> Dacarbazine
# This mirrors the existing "1000x" dilution columns (K/L/M) already present
# for the other compounds in the workbook.
# ---------------------------------------------------------------------------

# "condition" sheet (plain drug names) -- written first so the new unique
# shared strings "Temozolomide"/"Dacarbazine" are minted before the combined
# "1000x\n<drug>\n 1" labels used on the "merged" sheet.
$wsCondition = $wb.Worksheets.Item("condition")
$wsCondition.Range("K2").Value = "Temozolomide"
$wsCondition.Range("K3").Value = "Dacarbazine"

# Highlight the newly-added column on the condition sheet the way the author
# did: K2 gets a left/top/bottom accent-colored border, K3 gets a light blue
# fill -- a visual cue calling out the freshly entered values.
$k2Top = $wsCondition.Range("K2").Borders.Item(8)
$k2Top.LineStyle = 1
$k2Top.Color = 14136213
$k2Bottom = $wsCondition.Range("K2").Borders.Item(9)
$k2Bottom.LineStyle = 1
$k2Bottom.Color = 14136213
$wsCondition.Range("K2").Borders.Item(7).LineStyle = 1

$wsCondition.Range("K3").Borders.Item(7).LineStyle = 1
$wsCondition.Range("K3").Borders.Item(8).LineStyle = 1
$wsCondition.Range("K3").Borders.Item(9).LineStyle = 1
$wsCondition.Range("K3").Borders.Item(10).LineStyle = 1
$wsCondition.Range("K3").Interior.Color = 15918812

# "merged" sheet (dilution + drug name + replicate, combined label)
$wsMerged = $wb.Worksheets.Item("merged")
$wsMerged.Range("K2").Value = "1000x" + [char]10 + "Temozolomide" + [char]10 + " 1"
$wsMerged.Range("K3").Value = "1000x" + [char]10 + "Dacarbazine" + [char]10 + " 1"

# "sample" sheet (dilution factor only)
$wsSample = $wb.Worksheets.Item("sample")
$wsSample.Range("K2").Value = "1000x"
$wsSample.Range("K3").Value = "1000x"

# "replicate" sheet (replicate number)
$wsReplicate = $wb.Worksheets.Item("replicate")
$wsReplicate.Range("K2").Value = 1
$wsReplicate.Range("K3").Value = 1

# ---------------------------------------------------------------------------
# Print setup tweak on the "merged" sheet: fit to one page, landscape, 64%.
# ---------------------------------------------------------------------------
$ps = $wsMerged.PageSetup
$ps.Zoom = 64
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.PaperSize = 9
$ps.Orientation = 2

# ---------------------------------------------------------------------------
# Restore each sheet's last-used selection the way the author left them.
# ---------------------------------------------------------------------------
$wsSample.Activate()
$wsSample.Range("K4").Select()

$wsCondition.Activate()
$wsCondition.Range("K3").Select()

$wsReplicate.Activate()
$wsReplicate.Range("K3").Select()

$wsMerged.Activate()
$wsMerged.Range("K5").Select()

Write-Output "edit complete"
